# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2, 3) {
    $ws.Cells.Item($r, 4).Value = -0.0442     # D: historical_growth_revenue_last_5_years
    $ws.Cells.Item($r, 6).ClearContents()     # F: expected_growth_eps_next_5_years (removed)

    $ws.Cells.Item($r, 7).Value = 0.2273560491697218     # G: ebitdard_margin
    $ws.Cells.Item($r, 8).Value = 0.2273560491697218     # H: ebitda_margin
    $ws.Cells.Item($r, 9).Value = -0.528338708275598     # I: operating_margin
    $ws.Cells.Item($r, 10).Value = -0.528338708275598    # J: after_tax_operating_margin
    $ws.Cells.Item($r, 11).Value = -886.1                # K: trailing_net_income
    $ws.Cells.Item($r, 12).Value = -0.4777334483502265   # L: net_margin

    $ws.Cells.Item($r, 21).Value = 363.3                 # U: cash
    $ws.Cells.Item($r, 22).Value = 0.9010416666666667    # V: cash_market_cap
    $ws.Cells.Item($r, 23).Value = -1.539972193256865    # W: roe
    $ws.Cells.Item($r, 24).Value = 0.5175970166611015    # X: cost_equity
    $ws.Cells.Item($r, 25).Value = -2.057569209917967    # Y: roe_cost_equity
    $ws.Cells.Item($r, 26).Value = 0.2515370782703055    # Z: sales_invested_capital
    $ws.Cells.Item($r, 27).Value = -0.1328967750167512   # AA: roic
    $ws.Cells.Item($r, 28).Value = 0.06482396144267341   # AB: cost_capital
    $ws.Cells.Item($r, 29).Value = -0.1977207364594246   # AC: roic_cost_capital
    $ws.Cells.Item($r, 30).Value = 5544.5                # AD: debt_total
    $ws.Cells.Item($r, 31).Value = 3.963180547895941     # AE: debt_leases
    $ws.Cells.Item($r, 32).Value = 5548.463180547896     # AF: debt_total_inc_leases
    $ws.Cells.Item($r, 33).Value = 5185.163180547896     # AG: net_debt
    $ws.Cells.Item($r, 34).Value = 0.9322542308311738    # AH: debt_market_capital
    $ws.Cells.Item($r, 35).Value = 0.8238542368607101    # AI: debt_book_capital
    $ws.Cells.Item($r, 36).Value = 0.9278500721994112    # AJ: net_debt_market_capital
    $ws.Cells.Item($r, 37).Value = 0.8138104284080651    # AK: net_debt_book_capital
    $ws.Cells.Item($r, 38).Value = 316.9                 # AL: interest_expenses
    $ws.Cells.Item($r, 39).Value = 302.5                 # AM: net_interest_expenses
    $ws.Cells.Item($r, 40).Value = -8.299280015567231    # AN: debt_ebitda
    $ws.Cells.Item($r, 41).Value = -3.094351530451247    # AO: ebit_interest_expenses
    $ws.Cells.Item($r, 42).Value = -7.761407009067755    # AP: net_debt_ebitda
    $ws.Cells.Item($r, 43).Value = -3.241652892561984    # AQ: ebit_net_interest_expenses
}
